$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("error message")

# Insert a new row before existing row 5 (shifts old rows 5-7 down to 6-8)
$ws.Rows.Item(5).Insert()

# New row 5: Chinese (CH) message for beas_qc_1001
$ws.Range("A5").Value = "CH  "
$ws.Range("B5").Value = "NULL"
$ws.Range("C5").Value = "beas_qc_1001"

# New row 9 (appended at the end): English (E) message for beas_qc_1001
$ws.Range("A9").Value = "E   "
$ws.Range("B9").Value = "NULL"
$ws.Range("C9").Value = "beas_qc_1001"

# Shared strings are minted in this order: beas_qc_1001, English msg, Chinese msg
$ws.Range("E9").Value = "The Serial Number is error, can't find in goods receipt serial number list(OSRN)"
$ws.Range("E5").Value = "序列号输入错误，不在对应的采购收货序列号(OSRN)列表中"

$ws.Range("E6").Select()
